# C2 2.2.8 Started and Initial writing added
# Populates the "2.2.7 Professional Skills" sheet with a title block and
# three yearly ("Year : 2019-20" / "2018-19" / "2017-18") event tables.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2.2.7 Professional Skills")

$xlCenter = -4108

# ---- Wipe the old B4:F10 layout before laying out the new content ----
$ws.Range("B4:F10").Clear()

# ---- Column widths -------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 7.7109375
$ws.Columns.Item(7).ColumnWidth = 30.85546875

# ---- Title block (rows 4-6) -----------------------------------------
# Row 4: institute name, merged C4:G4
$r4 = $ws.Range("C4:G4")
$r4.ClearFormats()
$r4.Font.Name = "Verdana"
$r4.Font.Bold = $true
$r4.Font.Size = 16
$r4.HorizontalAlignment = $xlCenter
$r4.Merge()
$ws.Range("C4").Value = "Government Residential Womens Polytechinc, Latur"
$ws.Rows.Item(4).RowHeight = 19.5

# Row 5: blank spacer styled the same as row 4, NOT merged
$r5 = $ws.Range("C5:G5")
$r5.ClearFormats()
$r5.Font.Name = "Verdana"
$r5.Font.Bold = $true
$r5.Font.Size = 16
$r5.HorizontalAlignment = $xlCenter
$ws.Rows.Item(5).RowHeight = 19.5

# Row 6: format description, merged C6:G6, partially rich text
$r6 = $ws.Range("C6:G6")
$r6.ClearFormats()
$r6.Font.Name = "Verdana"
$r6.Font.Bold = $true
$r6.Font.Size = 12
$r6.HorizontalAlignment = $xlCenter
$r6.Merge()

$part1 = "Format : List of information about "
$part2 = "Events/Activity"
$part3 = " conducted for Womens/Students"
$ws.Range("C6").Value = $part1 + $part2 + $part3

$c6 = $ws.Range("C6")
$run2 = $c6.Characters($part1.Length + 1, $part2.Length)
$run2.Font.Name = "Verdana"
$run2.Font.Bold = $true
$run2.Font.Size = 12
$run2.Font.Underline = $true

$run3 = $c6.Characters($part1.Length + $part2.Length + 1, $part3.Length)
$run3.Font.Name = "Verdana"
$run3.Font.Bold = $true
$run3.Font.Size = 12

$ws.Rows.Item(6).RowHeight = 15.75

# ---- Helper data for the three yearly blocks -------------------------
$years = @(
    @{ Label = "Year : 2019-20"; LabelRow = 9;  HeaderRow = 10; FirstDataRow = 11 },
    @{ Label = "Year : 2018-19"; LabelRow = 19; HeaderRow = 20; FirstDataRow = 21 },
    @{ Label = "Year : 2017-18"; LabelRow = 29; HeaderRow = 30; FirstDataRow = 31 }
)

$headers = @("Sr.No ", "Name of the Event", "Date/Duration", "Conducting authority", "No. of beneficiaries", "Remark/ any additional info")

foreach ($year in $years) {
    $labelRow = $year.LabelRow
    $headerRow = $year.HeaderRow
    $firstDataRow = $year.FirstDataRow

    # Year label, merged B:D
    $lbl = $ws.Range("B" + $labelRow + ":D" + $labelRow)
    $lbl.ClearFormats()
    $lbl.Font.Name = "Verdana"
    $lbl.Font.Bold = $true
    $lbl.Font.Size = 11
    $lbl.Borders.LineStyle = 1
    $lbl.HorizontalAlignment = $xlCenter
    $lbl.Merge()
    $ws.Range("B" + $labelRow).Value = $year.Label

    # Header row
    for ($col = 2; $col -le 7; $col++) {
        $cell = $ws.Cells.Item($headerRow, $col)
        $cell.ClearFormats()
        $cell.Font.Name = "Verdana"
        $cell.Font.Bold = $true
        $cell.Font.Size = 10
        $cell.Borders.LineStyle = 1
        $cell.HorizontalAlignment = $xlCenter
        $cell.VerticalAlignment = $xlCenter
        $cell.Value = $headers[$col - 2]
    }

    # Data rows (6 rows, Sr.No. 1-6)
    for ($i = 0; $i -lt 6; $i++) {
        $row = $firstDataRow + $i
        for ($col = 2; $col -le 6; $col++) {
            $cell = $ws.Cells.Item($row, $col)
            $cell.ClearFormats()
            $cell.Borders.LineStyle = 1
            $cell.HorizontalAlignment = $xlCenter
            $cell.VerticalAlignment = $xlCenter
        }
        $ws.Cells.Item($row, 2).Value = $i + 1

        $remark = $ws.Cells.Item($row, 7)
        $remark.ClearFormats()
        $remark.Borders.LineStyle = 1
    }
}

# ---- Selection / active cell -----------------------------------------
$ws.Range("B9:D9").Select()

Write-Host "2.2.7 Professional Skills sheet populated"
